# The deck's presentation theme (ppt/theme/theme1.xml) currently carries the
# custom "Integral" look (clrScheme name="Integral", non-default accent
# colours). The commit swaps it back to the default PowerPoint "Office
# Theme" colour scheme (clrScheme name="Office": dk2=44546A, lt2=E7E6E6,
# accent1=5B9BD5, accent2=ED7D31, accent3=A5A5A5, accent4=FFC000,
# accent5=4472C4, accent6=70AD47, hlink=0563C1, folHlink=954F72 - dk1/lt1
# stay black/white). Font scheme and format scheme are already identical
# between the two themes, so only the twelve theme colours need updating.
#
# PowerPoint exposes the live 12-slot theme colour scheme through
# Design/ThemeColorScheme (mirrors a:clrScheme dk1,lt1,dk2,lt2,accent1-6,
# hlink,folHlink in that order). RGB is the usual VBA-style integer
# (R + G*256 + B*65536), so each target hex colour is converted below.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

# Target "Office Theme" palette, in clrScheme order (dk1, lt1, dk2, lt2,
# accent1..accent6, hlink, folHlink).
$officeThemeRgb = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRgb[$i - 1]
}
